$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray "NA" text in C68 so it becomes blank, matching the source export
$ws.Cells.Item(68, 3).Value2 = ""

# New rows of scraped data appended by the script run on 2025-05-08
$newRows = @(
    @(69, "2025-05-08", "eaux de surface", 237, 2),
    @(70, "2025-05-08", "substances actives", 237, 1),
    @(71, "2025-05-08", "substance active", 237, 1),
    @(72, "2025-05-08", "eaux de surface", 238, 2),
    @(73, "2025-05-08", "substances actives", 238, 1),
    @(74, "2025-05-08", "substance active", 238, 1),
    @(75, "2025-05-08", "eaux de surface", 239, 2),
    @(76, "2025-05-08", "substances actives", 239, 1),
    @(77, "2025-05-08", "substance active", 240, 1),
    @(78, "2025-05-08", "eaux de surface", 241, 2),
    @(79, "2025-05-08", "substances actives", 241, 1),
    @(80, "2025-05-08", "substance active", 241, 1),
    @(81, "2025-05-08", "eaux de surface", 242, 2),
    @(82, "2025-05-08", "substances actives", 242, 1),
    @(83, "2025-05-08", "substance active", 242, 1),
    @(84, "2025-05-08", "eaux de surface", 243, 2),
    @(85, "2025-05-08", "substances actives", 243, 1),
    @(86, "2025-05-08", "substance active", 243, 1),
    @(87, "2025-05-08", "eaux de surface", 244, 1),
    @(88, "2025-05-08", "substances actives", 245, 1),
    @(89, "2025-05-08", "eaux de surface", 245, 1),
    @(90, "2025-05-08", "herbicides", 245, 2),
    @(91, "2025-05-08", "eaux de surface", 247, 2),
    @(92, "2025-05-08", "substances actives", 247, 1),
    @(93, "2025-05-08", "substance active", 247, 1),
    @(94, "2025-05-08", "eaux de surface", 249, 2),
    @(95, "2025-05-08", "substances actives", 249, 1),
    @(96, "2025-05-08", "substance active", 249, 1),
    @(97, "2025-05-08", "eaux de surface", 250, 2),
    @(98, "2025-05-08", "substances actives", 250, 1),
    @(99, "2025-05-08", "substance active", 250, 1),
    @(100, "2025-05-08", "eaux de surface", 251, 2),
    @(101, "2025-05-08", "substances actives", 251, 1),
    @(102, "2025-05-08", "substance active", 251, 1),
    @(103, "2025-05-08", "eaux de surface", 252, 2),
    @(104, "2025-05-08", "substances actives", 252, 1),
    @(105, "2025-05-08", "substance active", 252, 1)
)

foreach ($row in $newRows) {
    $r = $row[0]

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value2 = $row[1]
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
}
